$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the current header text (B1:E1) and data block (B2:E5) first ---
$headers = @($ws.Range("B1").Value2, $ws.Range("C1").Value2, $ws.Range("D1").Value2, $ws.Range("E1").Value2)

$data = @(@(0,0,0,0), @(0,0,0,0), @(0,0,0,0), @(0,0,0,0))
for ($r = 0; $r -lt 4; $r++) {
    $row = 2 + $r
    $data[$r][0] = $ws.Cells.Item($row, 2).Value2  # old B
    $data[$r][1] = $ws.Cells.Item($row, 3).Value2  # old C
    $data[$r][2] = $ws.Cells.Item($row, 4).Value2  # old D
    $data[$r][3] = $ws.Cells.Item($row, 5).Value2  # old E
}

# --- The index column (A2:A5) loses its bordered/bold look - it's now plain data ---
$ws.Range("A2:A5").ClearFormats()

# --- Drop the now-unused column E entirely (values + formatting) ---
$ws.Range("E1:E5").Clear()

# --- Shift the header labels from B1:E1 into A1:D1 ---
$ws.Range("A1").Value = $headers[0]
$ws.Range("B1").Value = $headers[1]
$ws.Range("C1").Value = $headers[2]
$ws.Range("D1").Value = $headers[3]

# --- Shift the data block from B2:E5 into A2:D5 ---
for ($r = 0; $r -lt 4; $r++) {
    $row = 2 + $r
    $ws.Cells.Item($row, 1).Value = $data[$r][0]
    $ws.Cells.Item($row, 2).Value = $data[$r][1]
    $ws.Cells.Item($row, 3).Value = $data[$r][2]
    $ws.Cells.Item($row, 4).Value = $data[$r][3]
}

# --- A1 now also needs the bold/bordered/centered header look that B1:D1 already have ---
$a1 = $ws.Range("A1")
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Font.Bold = $true

$ws.Range("A1:D5").Select() | Out-Null
